$wb = $excel.ActiveWorkbook

# --- Sheet "URL" ---
$wsUrl = $wb.Worksheets.Item("URL")

$wsUrl.Range("A2").Value = "techable"
$wsUrl.Range("B2").Value = 0
$wsUrl.Range("E2").Value = 2
$wsUrl.Range("F2").Value = 2
$wsUrl.Range("G2").Value = "https://techable.jp/page/(pagenum)/?s=(keyword)"
$wsUrl.Range("H2").Value = "body > div> div > div > div > div > div > div > div > div > div> section > a"
$wsUrl.Range("I2").Value = "body > div> div > div > div > div > div > div > div > div > div> section > a>div>div>h3"
$wsUrl.Range("J2").Value = "body > div > div > div> div > div > div> time"
$wsUrl.Range("L2").Value = "body > div > div  div > div > div.te-cms-body"
$wsUrl.Range("O2").Value = "https://techable.jp"
$wsUrl.Range("Q2").Value = 14
$wsUrl.Range("U2").Value = "body > div.te-viewport > div.te-contents > div.te-layout > div.te-layout__col.te-layout__col--main > div > div.paging > div > a"
$wsUrl.Range("V2").Value = "14件"

$wsUrl.Range("C23").Select()

# --- Sheet "keyword" ---
$wsKeyword = $wb.Worksheets.Item("keyword")

$keywords = @(
  "ハプティクス",
  "ロボティクス",
  "ロボット",
  "触覚",
  "遠隔操作",
  "遠隔会議",
  "遠隔医療",
  "宇宙",
  "医療",
  "MR",
  "VR",
  "デジタルサイネージ",
  "建設",
  "ANA AVATAR",
  "シェアリングエコノミー",
  "VR広告"
)

$row = 4
foreach ($kw in $keywords) {
  $wsKeyword.Cells.Item($row, 1).Value = $kw
  $row = $row + 1
}

$wsKeyword.Range("A20").Select()

$wsUrl.Activate()
$wsUrl.Range("C23").Select()
